$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Create the new "Scan-Based Contrastive" worksheet by copying the
#    existing "Contrastive Pre-Trained" sheet (same layout/styles),
#    placing it after the last sheet, then renaming it.
# ------------------------------------------------------------------
$src = $wb.Worksheets.Item("Contrastive Pre-Trained")
$src.Copy($null, $src)
$new = $wb.Worksheets.Item($wb.Worksheets.Count)
$new.Name = "Scan-Based Contrastive"

# Populate the new sheet with the newly-collected scan-based
# contrastive performance numbers (2 additional trained models).
$new.Range("B2").Value = 0.9022
$new.Range("C2").Value = 8.1250999999999998

$new.Range("B3").Value = 0.89980000000000004
$new.Range("C3").Value = 10.0106

$new.Range("B4").Value = 0.89239999999999997
$new.Range("C4").Value = 13.8286

# Iteration 4 scan failed / no result recorded.
$new.Range("B5").ClearContents()
$new.Range("C5").ClearContents()

$new.Range("B6").Value = 0.82889999999999997
$new.Range("C6").Value = 13.5275

$new.Range("B7").Value = 0.86209999999999998
$new.Range("C7").Value = 14.0457

$new.Range("B8").Value = 0.88560000000000005
$new.Range("C8").Value = 9.8329000000000004

$new.Range("B9").Value = 0.86099999999999999
$new.Range("C9").Value = 9.8215000000000003

$new.Range("B10").Value = 0.87329999999999997
$new.Range("C10").Value = 12.109

$new.Range("B11").Value = 0.90800000000000003
$new.Range("C11").Value = 12.489000000000001

$new.Range("C12").Select()

# The source sheet ends up with its whole data table selected
# (e.g. selected via the table's "select all" corner) rather than
# being the active tab any more.
$src.Range("A1:C13").Select()

# ------------------------------------------------------------------
# 2. Update the "T-Tests" sheet with a new comparison block for
#    "Scan-Based Contrastive vs Standard".
# ------------------------------------------------------------------
$tt = $wb.Worksheets.Item("T-Tests")

$tt.Range("F9").Value = "Scan-Based Contrastive vs Standard"

$tt.Range("F10").Value = "Dice:"
$tt.Range("H10").Formula = "=_xlfn.T.TEST('Standard Training'!B2:B11, 'Scan-Based Contrastive'!B2:B11, 2, 3)"

$tt.Range("F11").Value = "Hausdorff:"
$tt.Range("H11").Formula = "=_xlfn.T.TEST('Standard Training'!C2:C11, 'Scan-Based Contrastive'!C2:C11, 2, 3)"

# ------------------------------------------------------------------
# 3. Restore focus to the "T-Tests" sheet (first tab) as the active
#    sheet/selection, matching the state the workbook was saved in.
# ------------------------------------------------------------------
$tt.Activate()
$tt.Range("H12").Select()

$wb.Application.CalculateFull()
